$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: convert from "Pending" (not-yet-recorded) to "Recorded" ---
# Copy formatting from row 10 (a "Recorded" row) onto row 11, then set its values.
$ws.Range("A10:I10").Copy()
$ws.Range("A11:I11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("G11").Value = "Sarah.Mahdy@med.asu.edu.eg"
$ws.Range("H11").Value = "2/216"
$ws.Range("I11").Value = "Recorded"

# --- Remaining cell value updates (recorder-email reorderings & recomputed stats) ---
$ws.Range("L6").Value = 67
$ws.Range("G8").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("L8").Value = 84
$ws.Range("L9").Value = "38.1%"
$ws.Range("G10").Value = "amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
$ws.Range("L10").Value = "35.8%"
$ws.Range("G15").Value = "hana.amr@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg"
$ws.Range("O15").Value = 10
$ws.Range("Q15").Value = 9
$ws.Range("R15").Value = "45.5%"
$ws.Range("S15").Value = "34.3%"
$ws.Range("G18").Value = "yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G19").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G20").Value = "yassmen.ahmed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Range("G21").Value = "yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G25").Value = "backup@backdoor.com, manar.montaser@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, gehanadel@med.asu.edu.eg"
$ws.Range("G30").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G32").Value = "amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
$ws.Range("G37").Value = "hana.amr@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg"
$ws.Range("G40").Value = "yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G41").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G42").Value = "yassmen.ahmed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Range("G46").Value = "nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G54").Value = "amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G58").Value = "Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G62").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G63").Value = "ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G68").Value = "nourhan.mahmoud@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("G76").Value = "amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G80").Value = "Amr-Saeed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G81").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg, enas.omran@med.asu.edu.eg"
$ws.Range("G84").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
$ws.Range("G85").Value = "ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G90").Value = "Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg"
$ws.Range("G96").Value = "mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G98").Value = "amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G106").Value = "wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G107").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G108").Value = "yassmen.ahmed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Range("G112").Value = "Mohammedeltanany@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg"
$ws.Range("G118").Value = "mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G120").Value = "amany.raafat@med.asu.edu.eg, merna.said@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G128").Value = "wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G129").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G130").Value = "yassmen.ahmed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg, marinasorial@med.asu.edu.eg"
$ws.Range("G134").Value = "Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G140").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G142").Value = "merna.said@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G150").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G151").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G152").Value = "marina_atef@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg"
$ws.Range("G156").Value = "Mohammedeltanany@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G162").Value = "aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg"
$ws.Range("G164").Value = "merna.said@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg"
$ws.Range("G172").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G173").Value = "yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G174").Value = "marina_atef@med.asu.edu.eg, Wafaa.ebida@med.asu.edu.eg"
